$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item(1).Name = "HT1"
$wb.Worksheets.Item(2).Name = "HT2"
